# Insert a new data row at row 50 (pushing the existing rows 50-91 down to 51-92)
# and populate it with a new "Acelga" price record for
# "Agrícola del Norte S.A. de Arica".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("50:50").Insert()

$ws.Cells.Item(50, 1).Value  = 1
$ws.Cells.Item(50, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(50, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(50, 4).Value  = 45049
$ws.Cells.Item(50, 5).Value  = 15
$ws.Cells.Item(50, 6).Value  = 100112009
$ws.Cells.Item(50, 7).Value  = "Acelga"
$ws.Cells.Item(50, 8).Value  = "Sin especificar"
$ws.Cells.Item(50, 9).Value  = "Segunda"
$ws.Cells.Item(50, 10).Value = 400
$ws.Cells.Item(50, 11).Value = 3000
$ws.Cells.Item(50, 12).Value = 3500
$ws.Cells.Item(50, 13).Value = 3312
$ws.Cells.Item(50, 14).Value = "$/atado 2,5 a 3 kilos"
$ws.Cells.Item(50, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(50, 16).Value = 1104
$ws.Cells.Item(50, 17).Value = 3
$ws.Cells.Item(50, 18).Value = "Hortaliza"
